$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the two new sheets ('ODI Batting' data-holder + 'ODI
# Batting Extra') appended AFTER the existing sheet, so sheetId is
# assigned sequentially (1,2,3) matching the final tab order. The
# original sheet (sheetId=1) keeps its identity and becomes 'Player
# Info'; the freshly appended sheets take sheetId 2 and 3.
# ---------------------------------------------------------------------
$origSheet = $wb.Worksheets.Item("ODI Batting")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingSheet = $wb.Worksheets.Add($null, $lastSheet)
$battingSheet.Name = "ODI Batting NEW"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$extraSheet = $wb.Worksheets.Add($null, $lastSheet2)
$extraSheet.Name = "ODI Batting Extra"

# ---------------------------------------------------------------------
# Step 2: populate the new 'ODI Batting NEW' sheet with the original
# batting data, transformed: header D1 renamed MATCH_CARD_LINK ->
# MATCH_CODE, column D values reduced from the full scorecard URL down
# to the bare numeric match code, and the handful of originally-blank
# INNING_NUMBER cells left untouched (absent).
# ---------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting NEW")

# Header row (bold / bordered / centered -- copy exact style from the
# original header row so the new header matches formatting).
$battingHeader = @('MATCH_NUMBER', 'INNING_NUMBER', 'MATCH_DATE', 'MATCH_CODE', 'MATCH_INNING', 'OPPONENT', 'VENUE', 'DISMISSAL', 'RUNS_SCORED', 'BALLS_FACED')
for ($c = 0; $c -lt $battingHeader.Length; $c++) {
    $srcCell = $origSheet.Cells.Item(1, $c + 1)
    $dstCell = $battingSheet.Cells.Item(1, $c + 1)
    $srcCell.Copy($dstCell)
    $dstCell.Value = $battingHeader[$c]
}

# Force the whole data range to text format first (so numeric-looking
# values like match numbers / scores are preserved as text, matching
# the source workbook), then fill in the 122 data rows.
$battingSheet.Range("A2:J123").NumberFormat = "@"

$battingData = @(
    @('1', '1', '26/06/2009', '2975', '2nd', 'India', 'Sabina Park', 'c R P Singh b Harbhajan Singh', '19', '16'),
    @('2', $null, '28/06/2009', '2976', '2nd', 'India', 'Sabina Park', 'did not bat', '-', '-'),
    @('3', '2', '03/07/2009', '2977', '1st', 'India', 'Daren Sammy National Cricket Stadium', 'b A Nehra', '21', '15'),
    @('4', $null, '05/07/2009', '2978', '1st', 'India', 'Daren Sammy National Cricket Stadium', 'did not bat', '-', '-'),
    @('5', '3', '14/03/2010', '3106', '2nd', 'Zimbabwe', 'Arnos Vale Ground', 'st T Taibu b G A Lamb', '14', '35'),
    @('6', '4', '13/04/2010', '3107', '1st', 'Canada', 'Sabina Park', 'b Rizwan Cheema', '74', '75'),
    @('7', '5', '15/04/2010', '3108', '2nd', 'Ireland', 'Sabina Park', 'not out', '28*', '36'),
    @('8', '6', '28/05/2010', '3112', '2nd', 'South Africa', 'Windsor Park', 'lbw b M Morkel', '1', '6'),
    @('9', '7', '30/05/2010', '3114', '1st', 'South Africa', 'Windsor Park', 'not out', '45*', '31'),
    @('10', '8', '03/06/2010', '3117', '1st', 'South Africa', 'Queen''s Park Oval', 'c M V Boucher b R E van der Merwe', '17', '47'),
    @('11', '9', '31/01/2011', '3227', '1st', 'Sri Lanka', 'Sinhalese Sports Club Ground', 'run out', '8', '22'),
    @('12', '10', '03/02/2011', '3231', '1st', 'Sri Lanka', 'Sinhalese Sports Club Ground', 'lbw b S L Malinga', '39', '64'),
    @('13', '11', '06/02/2011', '3234', '2nd', 'Sri Lanka', 'Sinhalese Sports Club Ground', 'lbw b N L T C Perera', '79', '91'),
    @('14', '12', '24/02/2011', '3241', '1st', 'South Africa', 'Arun Jaitley Stadium', 'lbw b J Botha', '73', '82'),
    @('15', '13', '28/02/2011', '3247', '1st', 'Netherlands', 'Arun Jaitley Stadium', 'c A N Kervezee b P M Seelaar', '30', '38'),
    @('16', '14', '04/03/2011', '3253', '2nd', 'Bangladesh', 'Shere Bangla National Stadium', 'not out', '9*', '26'),
    @('17', '15', '11/03/2011', '3261', '1st', 'Ireland', 'Punjab Cricket Association IS Bindra Stadium', 'b K J O''Brien', '0', '3'),
    @('18', '16', '17/03/2011', '3270', '2nd', 'England', 'MA Chidambaram Stadium', 'c A J Strauss b J C Tredwell', '5', '15'),
    @('19', '17', '20/03/2011', '3276', '2nd', 'India', 'MA Chidambaram Stadium', 'c Harbhajan Singh b S K Raina', '22', '29'),
    @('20', '18', '23/03/2011', '3277', '1st', 'Pakistan', 'Shere Bangla National Stadium', 'lbw b Mohammad Hafeez', '0', '3'),
    @('21', '19', '23/04/2011', '3287', '1st', 'Pakistan', 'Daren Sammy National Cricket Stadium', 'run out', '67', '109'),
    @('22', '20', '25/04/2011', '3288', '1st', 'Pakistan', 'Daren Sammy National Cricket Stadium', 'c Misbah-ul-Haq b Saeed Ajmal', '12', '18'),
    @('23', '21', '28/04/2011', '3289', '1st', 'Pakistan', 'Kensington Oval', 'lbw b Mohammad Hafeez', '47', '77'),
    @('24', '22', '02/05/2011', '3290', '2nd', 'Pakistan', 'Kensington Oval', 'c Tanvir Ahmed b Junaid Khan', '21', '26'),
    @('25', $null, '05/05/2011', '3291', '2nd', 'Pakistan', 'Providence Stadium', 'did not bat', '-', '-'),
    @('26', '23', '06/06/2011', '3294', '1st', 'India', 'Queen''s Park Oval', 'c R G Sharma b M M Patel', '4', '8'),
    @('27', '24', '13/06/2011', '3297', '1st', 'India', 'Sir Vivian Richards Stadium', 'c M K Tiwary b A Mishra', '15', '36'),
    @('28', '25', '16/06/2011', '3298', '2nd', 'India', 'Sabina Park', 'st P A Patel b A Mishra', '86', '99'),
    @('29', '26', '13/10/2011', '3334', '1st', 'Bangladesh', 'Shere Bangla National Stadium', 'b Shafiul Islam', '20', '19'),
    @('30', '27', '15/10/2011', '3336', '2nd', 'Bangladesh', 'Shere Bangla National Stadium', 'not out', '1*', '12'),
    @('31', '28', '18/10/2011', '3338', '1st', 'Bangladesh', 'Zahur Ahmed Chowdhury Stadium', 'c †Mushfiqur Rahim b Shafiul Islam', '0', '3'),
    @('32', '29', '29/11/2011', '3353', '1st', 'India', 'Barabati Stadium', 'b S K Raina', '60', '74'),
    @('33', '30', '02/12/2011', '3355', '1st', 'India', 'Dr YS Rajasekhara Reddy Cricket Stadium', 'c R Ashwin b R Vinay Kumar', '13', '17'),
    @('34', '31', '05/12/2011', '3357', '1st', 'India', 'Narendra Modi Stadium', 'retired hurt', '26*', '41'),
    @('35', '32', '16/03/2012', '3398', '2nd', 'Australia', 'Arnos Vale Ground', 'run out', '4', '15'),
    @('36', '33', '18/03/2012', '3400', '2nd', 'Australia', 'Arnos Vale Ground', 'b X J Doherty', '16', '28'),
    @('37', '34', '20/03/2012', '3402', '2nd', 'Australia', 'Arnos Vale Ground', 'c G J Bailey b S R Watson', '0', '3'),
    @('38', '35', '23/03/2012', '3404', '1st', 'Australia', 'Daren Sammy National Cricket Stadium', 'c D J Hussey b S R Watson', '25', '50'),
    @('39', '36', '25/03/2012', '3405', '2nd', 'Australia', 'Daren Sammy National Cricket Stadium', 'c †M S Wade b C J McKay', '3', '21'),
    @('40', '37', '16/06/2012', '3412', '2nd', 'England', 'The Rose Bowl', 'not out', '8*', '13'),
    @('41', '38', '30/11/2012', '3450', '1st', 'Bangladesh', 'Sheikh Abu Naser Stadium', 'run out', '35', '41'),
    @('42', '39', '02/12/2012', '3451', '2nd', 'Bangladesh', 'Sheikh Abu Naser Stadium', 'c †Mushfiqur Rahim b Abdur Razzak', '28', '41'),
    @('43', '40', '05/12/2012', '3452', '2nd', 'Bangladesh', 'Shere Bangla National Stadium', 'c Sohag Gazi b Naeem Islam', '13', '23'),
    @('44', '41', '07/12/2012', '3453', '1st', 'Bangladesh', 'Shere Bangla National Stadium', 'c & b Sohag Gazi', '34', '61'),
    @('45', '42', '08/12/2012', '3454', '1st', 'Bangladesh', 'Shere Bangla National Stadium', 'c Mominul Haque b Mahmudullah', '51', '108'),
    @('46', '43', '01/02/2013', '3471', '1st', 'Australia', 'WACA Ground', 'c M J Clarke b J P Faulkner', '11', '25'),
    @('47', '44', '03/02/2013', '3472', '2nd', 'Australia', 'WACA Ground', 'lbw b M A Starc', '14', '20'),
    @('48', '45', '06/02/2013', '3473', '2nd', 'Australia', 'Manuka Oval', 'b J P Faulkner', '86', '96'),
    @('49', '46', '08/02/2013', '3474', '1st', 'Australia', 'Sydney Cricket Ground', 'c A J Finch b M G Johnson', '4', '16'),
    @('50', '47', '10/02/2013', '3475', '2nd', 'Australia', 'Melbourne Cricket Ground', 'c J P Faulkner b X J Doherty', '33', '50'),
    @('51', '48', '22/02/2013', '3478', '1st', 'Zimbabwe', 'National Cricket Stadium (Grenada)', 'not out', '100*', '71'),
    @('52', '49', '24/02/2013', '3480', '2nd', 'Zimbabwe', 'National Cricket Stadium (Grenada)', 'c †B R M Taylor b H Masakadza', '0', '2'),
    @('53', '50', '26/02/2013', '3481', '2nd', 'Zimbabwe', 'National Cricket Stadium (Grenada)', 'not out', '72*', '103'),
    @('54', '51', '07/06/2013', '3506', '2nd', 'Pakistan', 'Kennington Oval', 'c †Kamran Akmal b Mohammad Irfan', '0', '4'),
    @('55', '52', '11/06/2013', '3510', '1st', 'India', 'Kennington Oval', 'st M S Dhoni b R Ashwin', '35', '83'),
    @('56', '53', '14/06/2013', '3513', '2nd', 'South Africa', 'Sophia Gardens', 'run out', '12', '7'),
    @('57', '54', '28/06/2013', '3520', '2nd', 'Sri Lanka', 'Sabina Park', 'run out', '27', '30'),
    @('58', '55', '30/06/2013', '3521', '2nd', 'India', 'Sabina Park', 'c S Dhawan b R Ashwin', '55', '78'),
    @('59', '56', '05/07/2013', '3525', '2nd', 'India', 'Queen''s Park Oval', 'c R Ashwin b B Kumar', '1', '7'),
    @('60', '57', '07/07/2013', '3527', '2nd', 'Sri Lanka', 'Queen''s Park Oval', 'c B M A J Mendis b A D Mathews', '70', '84'),
    @('61', '58', '14/07/2013', '3531', '2nd', 'Pakistan', 'Providence Stadium', 'c †Umar Akmal b Mohammad Irfan', '5', '18'),
    @('62', '59', '16/07/2013', '3532', '1st', 'Pakistan', 'Providence Stadium', 'b Saeed Ajmal', '54', '81'),
    @('63', '60', '19/07/2013', '3533', '2nd', 'Pakistan', 'Daren Sammy National Cricket Stadium', 'c Mohammad Hafeez b Wahab Riaz', '17', '40'),
    @('64', '61', '21/07/2013', '3535', '1st', 'Pakistan', 'Daren Sammy National Cricket Stadium', 'c Saeed Ajmal b Wahab Riaz', '9', '18'),
    @('65', '62', '24/07/2013', '3538', '1st', 'Pakistan', 'Daren Sammy National Cricket Stadium', 'c †Umar Akmal b Junaid Khan', '9', '21'),
    @('66', '63', '21/11/2013', '3580', '1st', 'India', 'Nehru Stadium (Kochi)', 'b Mohammed Shami', '59', '77'),
    @('67', '64', '24/11/2013', '3581', '2nd', 'India', 'Dr YS Rajasekhara Reddy Cricket Stadium', 'c †M S Dhoni b R Ashwin', '50', '54'),
    @('68', '65', '27/11/2013', '3583', '1st', 'India', 'Green Park', 'not out', '51*', '53'),
    @('69', '66', '26/12/2013', '3593', '2nd', 'New Zealand', 'Eden Park', 'c K S Williamson b M J McClenaghan', '14', '18'),
    @('70', '67', '23/02/2014', '3617', '2nd', 'Ireland', 'Sabina Park', 'c †G C Wilson b S R Thompson', '1', '9'),
    @('71', '68', '28/02/2014', '3622', '1st', 'England', 'Sir Vivian Richards Stadium', 'lbw b M M Ali', '2', '12'),
    @('72', '69', '02/03/2014', '3625', '1st', 'England', 'Sir Vivian Richards Stadium', 'b S C J Broad', '13', '34'),
    @('73', '70', '05/03/2014', '3629', '2nd', 'England', 'Sir Vivian Richards Stadium', 'c R S Bopara b M M Ali', '16', '20'),
    @('74', '71', '20/08/2014', '3655', '2nd', 'Bangladesh', 'National Cricket Stadium (Grenada)', 'c †Mushfiqur Rahim b Al-Amin Hossain', '7', '14'),
    @('75', '72', '22/08/2014', '3657', '1st', 'Bangladesh', 'National Cricket Stadium (Grenada)', 'lbw b Sohag Gazi', '53', '82'),
    @('76', '73', '25/08/2014', '3661', '1st', 'Bangladesh', 'Warner Park', 'c Imrul Kayes b Mahmudullah', '124', '127'),
    @('77', '74', '08/10/2014', '3678', '1st', 'India', 'Nehru Stadium (Kochi)', 'c S Dhawan b A Mishra', '28', '45'),
    @('78', '75', '11/10/2014', '3680', '2nd', 'India', 'Arun Jaitley Stadium', 'b Mohammed Shami', '26', '44'),
    @('79', '76', '17/10/2014', '3683', '2nd', 'India', 'Himachal Pradesh Cricket Association Stadium', 'b A R Patel', '40', '51'),
    @('80', '77', '16/02/2015', '3752', '1st', 'Ireland', 'Saxton Oval', 'run out', '0', '0'),
    @('81', '78', '21/02/2015', '3757', '1st', 'Pakistan', 'Hagley Oval', 'retired hurt', '49*', '78'),
    @('82', '79', '01/11/2015', '3852', '1st', 'Sri Lanka', 'R Premadasa Stadium', 'b B A W Mendis', '38', '58'),
    @('83', '80', '04/11/2015', '3853', '1st', 'Sri Lanka', 'R Premadasa Stadium', 'c T M Dilshan b T A M Siriwardana', '21', '33'),
    @('84', '81', '07/11/2015', '3855', '1st', 'Sri Lanka', 'Pallekele International Cricket Stadium', 'b R A S Lakmal', '4', '7'),
    @('85', '82', '03/06/2016', '3892', '2nd', 'South Africa', 'Providence Stadium', 'lbw b A M Phangiso', '30', '69'),
    @('86', '83', '05/06/2016', '3893', '1st', 'Australia', 'Providence Stadium', 'c A Zampa b M R Marsh', '19', '31'),
    @('87', '84', '13/06/2016', '3898', '2nd', 'Australia', 'Warner Park', 'c †M S Wade b A Zampa', '39', '63'),
    @('88', '85', '15/06/2016', '3900', '2nd', 'South Africa', 'Warner Park', 'b W D Parnell', '11', '19'),
    @('89', '86', '21/06/2016', '3905', '1st', 'Australia', 'Kensington Oval', 'c S P D Smith b J R Hazlewood', '15', '20'),
    @('90', '87', '24/06/2016', '3907', '1st', 'South Africa', 'Kensington Oval', 'c F du Plessis b C H Morris', '102', '103'),
    @('91', '88', '26/06/2016', '3909', '2nd', 'Australia', 'Kensington Oval', 'c †M S Wade b M R Marsh', '6', '20'),
    @('92', '89', '30/09/2016', '3939', '2nd', 'Pakistan', 'Sharjah Cricket Stadium', 'b Mohammad Nawaz', '12', '28'),
    @('93', '90', '02/10/2016', '3943', '2nd', 'Pakistan', 'Sharjah Cricket Stadium', 'run out', '61', '74'),
    @('94', '91', '05/10/2016', '3944', '2nd', 'Pakistan', 'Sheikh Zayed Stadium', 'c †Sarfaraz Ahmed b Wahab Riaz', '17', '40'),
    @('95', '92', '09/12/2018', '4228', '1st', 'Bangladesh', 'Shere Bangla National Stadium', 'c Tamim Iqbal b Mashrafe Mortaza', '19', '51'),
    @('96', '93', '11/12/2018', '4229', '2nd', 'Bangladesh', 'Shere Bangla National Stadium', 'b Rubel Hossain', '27', '43'),
    @('97', '94', '14/12/2018', '4230', '1st', 'Bangladesh', 'Sylhet Stadium', 'b Mehidy Hasan Miraz', '10', '26'),
    @('98', '95', '20/02/2019', '4253', '1st', 'England', 'Kensington Oval', 'c C R Woakes b A U Rashid', '40', '30'),
    @('99', '96', '22/02/2019', '4254', '1st', 'England', 'Kensington Oval', 'run out', '25', '37'),
    @('100', $null, '25/02/2019', '4255', '1st', 'England', 'National Cricket Stadium (Grenada)', 'did not bat', '-', '-'),
    @('101', '97', '27/02/2019', '4256', '2nd', 'England', 'National Cricket Stadium (Grenada)', 'c B A Stokes b M A Wood', '61', '59'),
    @('102', '98', '01/03/2019', '4260', '2nd', 'England', 'Daren Sammy National Cricket Stadium', 'not out', '7*', '14'),
    @('103', '99', '05/05/2019', '4285', '1st', 'Ireland', 'Castle Avenue', 'not out', '9*', '7'),
    @('104', '100', '07/05/2019', '4286', '1st', 'Bangladesh', 'Castle Avenue', 'c †Mushfiqur Rahim b Shakib Al Hasan', '1', '4'),
    @('105', '101', '11/05/2019', '4291', '2nd', 'Ireland', 'The Village', 'c sub b P R Stirling', '17', '25'),
    @('106', '102', '13/05/2019', '4293', '1st', 'Bangladesh', 'The Village', 'lbw b Mehidy Hasan Miraz', '6', '13'),
    @('107', '103', '17/05/2019', '4296', '1st', 'Bangladesh', 'The Village', 'not out', '3*', '3'),
    @('108', '104', '31/05/2019', '4304', '2nd', 'Pakistan', 'Trent Bridge', 'c Babar Azam b Mohammad Amir', '0', '4'),
    @('109', $null, '10/06/2019', '4317', ' ', 'South Africa', 'The Rose Bowl', 'did not bat', '-', '-'),
    @('110', '105', '17/06/2019', '4325', '1st', 'Bangladesh', 'The Cooper Associates County Ground', 'b Mohammad Saifuddin', '19', '15'),
    @('111', '106', '22/02/2020', '4413', '1st', 'Sri Lanka', 'Sinhalese Sports Club Ground', 'run out', '39', '52'),
    @('112', '107', '26/02/2020', '4414', '2nd', 'Sri Lanka', 'Mahinda Rajapaksha International Cricket Stadium', 'c D de Silva b P W H de Silva', '16', '23'),
    @('113', '108', '01/03/2020', '4417', '2nd', 'Sri Lanka', 'Pallekele International Cricket Stadium', 'run out', '8', '12'),
    @('114', '109', '10/03/2021', '4449', '2nd', 'Sri Lanka', 'Sir Vivian Richards Stadium', 'not out', '37*', '47'),
    @('115', '110', '12/03/2021', '4450', '2nd', 'Sri Lanka', 'Sir Vivian Richards Stadium', 'c K N A Bandara b N L T C Perera', '10', '10'),
    @('116', '111', '14/03/2021', '4451', '2nd', 'Sri Lanka', 'Sir Vivian Richards Stadium', 'c D Karunaratne b R A S Lakmal', '102', '132'),
    @('117', '112', '20/07/2021', '4483', '2nd', 'Australia', 'Kensington Oval', 'c M R Marsh b J R Hazlewood', '2', '9'),
    @('118', '113', '24/07/2021', '4484', '2nd', 'Australia', 'Kensington Oval', 'b M A Starc', '0', '1'),
    @('119', '114', '26/07/2021', '4486', '1st', 'Australia', 'Kensington Oval', 'c J R Philippe b A C Agar', '18', '36'),
    @('120', '115', '06/02/2022', '4533', '1st', 'India', 'Narendra Modi Stadium', 'lbw b Washington Sundar', '18', '34'),
    @('121', '116', '09/02/2022', '4535', '2nd', 'India', 'Narendra Modi Stadium', 'c †R R Pant b M Prasidh Krishna', '1', '3'),
    @('122', '117', '11/02/2022', '4536', '2nd', 'India', 'Narendra Modi Stadium', 'c V Kohli b M Prasidh Krishna', '19', '30')
)

for ($r = 0; $r -lt $battingData.Length; $r++) {
    $row = $battingData[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ($null -ne $val) {
            $battingSheet.Cells.Item($r + 2, $c + 1).Value = $val
        }
    }
}

# ---------------------------------------------------------------------
# Step 3: populate 'ODI Batting Extra' with the new per-innings extra
# stats. Column B (BATTING_POSITION) is numeric; the remaining columns
# are text (kept as '@' number format so values like '0'/'1' don't
# silently become numbers). Rows with no data keep their cells present
# but blank (mirrors the source, where every row has all 6 cells).
# ---------------------------------------------------------------------
$extraSheet = $wb.Worksheets.Item("ODI Batting Extra")

$extraHeader = @('MATCH_CODE', 'BATTING_POSITION', 'NUM_4', 'NUM_6', 'PERCENT_RUNS_OF_TOTAL', 'MAN_OF_MATCH')
for ($c = 0; $c -lt $extraHeader.Length; $c++) {
    $srcCell = $origSheet.Cells.Item(1, $c + 1)
    $dstCell = $extraSheet.Cells.Item(1, $c + 1)
    $srcCell.Copy($dstCell)
    $dstCell.Value = $extraHeader[$c]
}

$extraSheet.Range("A2:A21").NumberFormat = "@"
$extraSheet.Range("C2:F21").NumberFormat = "@"

$extraData = @(
    @('4285', '3', '1', '0', '2.36%', 'NO'),
    @('4286', '3', '0', '0', '0.38%', 'NO'),
    @('4291', '3', '1', '0', '5.14%', 'NO'),
    @('4293', '3', '1', '0', '2.43%', 'NO'),
    @('4296', '3', '0', '0', '1.97%', 'NO'),
    @('4304', '3', '0', '0', $null, 'NO'),
    @('4317', $null, $null, $null, $null, 'NO'),
    @('4325', $null, $null, $null, $null, 'NO'),
    @('4413', '3', '3', '2', '13.49%', 'NO'),
    @('4414', $null, $null, $null, $null, 'NO'),
    @('4417', '5', '0', '0', '2.66%', 'NO'),
    @('4449', '3', '2', '1', '15.68%', 'NO'),
    @('4450', $null, $null, $null, $null, 'NO'),
    @('4451', '4', '5', '4', '36.96%', 'YES'),
    @('4483', '4', '0', '0', '1.63%', 'NO'),
    @('4484', '3', '0', '0', $null, 'NO'),
    @('4486', $null, $null, $null, $null, 'NO'),
    @('4533', $null, $null, $null, $null, 'NO'),
    @('4535', '3', '0', '0', '0.52%', 'NO'),
    @('4536', '3', '3', '0', '11.24%', 'NO')
)

for ($r = 0; $r -lt $extraData.Length; $r++) {
    $row = $extraData[$r]
    $rowNum = $r + 2
    # Column A: MATCH_CODE (text)
    if ($null -ne $row[0]) { $extraSheet.Cells.Item($rowNum, 1).Value = $row[0] }
    # Column B: BATTING_POSITION (numeric)
    if ($null -ne $row[1]) {
        $extraSheet.Cells.Item($rowNum, 2).Value = [int]$row[1]
    } else {
        $extraSheet.Cells.Item($rowNum, 2).NumberFormat = "@"
        $extraSheet.Cells.Item($rowNum, 2).Value = ""
    }
    # Columns C-E: text, blank-but-present when missing
    for ($c = 2; $c -le 4; $c++) {
        if ($null -ne $row[$c]) {
            $extraSheet.Cells.Item($rowNum, $c + 1).Value = $row[$c]
        } else {
            $extraSheet.Cells.Item($rowNum, $c + 1).NumberFormat = "@"
            $extraSheet.Cells.Item($rowNum, $c + 1).Value = ""
        }
    }
    # Column F: MAN_OF_MATCH (text, always present)
    $extraSheet.Cells.Item($rowNum, 6).Value = $row[5]
}

# ---------------------------------------------------------------------
# Step 4: turn the original sheet (still sheetId 1) into 'Player Info':
# grab a copy of the real header style (A1 on the new batting sheet
# already carries the exact bold/border/center style) before clearing,
# then clear out the old batting data and write the new 4-column player
# record, then rename it last (renaming earlier would break the
# Worksheets.Item("ODI Batting") lookups used above).
# ---------------------------------------------------------------------
$headerStyleSrc = $battingSheet.Range("A1")
$origSheet.Cells.Clear()

$playerHeader = @('ID', 'NAME', 'BATTING_HAND', 'BOWL_STYLE')
for ($c = 0; $c -lt $playerHeader.Length; $c++) {
    $cell = $origSheet.Cells.Item(1, $c + 1)
    $headerStyleSrc.Copy($cell)
    $cell.Value = $playerHeader[$c]
}

$origSheet.Range("A2:D2").NumberFormat = "@"
$playerData = @('3703', 'Darren Michael Bravo', 'Left Handed', 'Right Arm Medium Fast')
for ($c = 0; $c -lt $playerData.Length; $c++) {
    $origSheet.Cells.Item(2, $c + 1).Value = $playerData[$c]
}

# ---------------------------------------------------------------------
# Step 5: final renames, now that all data has been written using the
# original/stable sheet names as lookup keys.
# ---------------------------------------------------------------------
$origSheet.Name = "Player Info"
$battingSheet.Name = "ODI Batting"

Write-Output "done"
